$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

function Set-TextCell($range, $text) {
    # Force the cell to Text format first so numeric-looking strings
    # (e.g. "605.91") are kept verbatim instead of being parsed as numbers.
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $text
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.732.68"
$ws.Range("E2").Value = "  +0.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.520.13"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextCell "D5" "605.91"
$ws.Range("E5").Value = "  -0.73%  "

# Row 6 - Solana
Set-TextCell "D6" "196.41"
$ws.Range("E6").Value = "  +5.69%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.39%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -7.18%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -0.22%  "

# Row 11 - Avalanche
Set-TextCell "D11" "53.73"
$ws.Range("E11").Value = "  +1.25%  "

# Row 12 - ShibaInu
Set-TextCell "D12" "0.0000301"
$ws.Range("E12").Value = "  -2.40%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  -0.07%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.082.37"
$ws.Range("E14").Value = "  +1.35%  "

# Row 15 - BitcoinCash
Set-TextCell "D15" "596.82"
$ws.Range("E15").Value = "  -0.83%  "

# Row 16 - Uniswap
$ws.Range("E16").Value = "  +1.52%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.956.84"
$ws.Range("E17").Value = "  +0.84%  "

# Row 18 - Chainlink
Set-TextCell "D18" "19.08"
$ws.Range("E18").Value = "  +1.22%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.533.22"
$ws.Range("E19").Value = "  +1.53%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +1.39%  "

# Row 21 - Polygon
Set-TextCell "D21" "0.992"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22 - InternetComputer(DFINITY)
Set-TextCell "D22" "18.31"
$ws.Range("E22").Value = "  +6.41%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  +5.13%  "

# Rows 24 & 25 - swap PancakeSwap/Litecoin with updated figures
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D24" "102.49"
$ws.Range("E24").Value = "  -3.26%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D25" "4.67"
$ws.Range("E25").Value = "  +0.65%  "

# Row 26 - ImmutableX
$ws.Range("E26").Value = "  +5.41%  "

# Row 27 - RenderToken
Set-TextCell "D27" "10.87"
$ws.Range("E27").Value = "  -0.68%  "

# Row 28 - Filecoin
Set-TextCell "D28" "9.60"
$ws.Range("E28").Value = "  -1.45%  "

# Row 29 - EthereumClassic
Set-TextCell "D29" "33.38"
$ws.Range("E29").Value = "  -0.45%  "

# Row 30 - dogwifhat
Set-TextCell "D30" "4.31"
$ws.Range("E30").Value = "  +10.84%  "

# Row 31 - NEARProtocol
Set-TextCell "D31" "7.08"
$ws.Range("E31").Value = "  +1.71%  "

# Row 32 - Cosmos
Set-TextCell "D32" "12.43"
$ws.Range("E32").Value = "  +0.18%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -0.10%  "

# Row 34 - OKB
Set-TextCell "D34" "63.10"
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - PEPE
$ws.Range("D35").Value = "0.0{0}0853" -f $sub3
$ws.Range("E35").Value = "  +10.00%  "

# Row 36 - Maker
$ws.Range("D36").Value = "3.739.70"
$ws.Range("E36").Value = "  +3.85%  "

# Row 37 - Fetch.AI
$ws.Range("E37").Value = "  -3.20%  "

# Row 38 - Dai
Set-TextCell "D38" "0.999"
$ws.Range("E38").Value = "  +0.15%  "

# Row 39 - Stacks
Set-TextCell "D39" "3.63"
$ws.Range("E39").Value = "  +0.33%  "

# Row 40 - TheGraph
Set-TextCell "D40" "0.392"
$ws.Range("E40").Value = "  -1.01%  "

# Row 41 - InjectiveProtocol
Set-TextCell "D41" "36.55"
$ws.Range("E41").Value = "  -0.54%  "

# Row 42 - Bittensor
Set-TextCell "D42" "487.97"
$ws.Range("E42").Value = "  -6.86%  "

# Row 43 - Kaspa
Set-TextCell "D43" "0.133"
$ws.Range("E43").Value = "  -3.00%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -1.59%  "

# Rows 45 & 46 - swap ThetaToken/Stellar with updated figures
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D45" "0.140"
$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell "D46" "2.82"
$ws.Range("E46").Value = "  -4.41%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  -1.68%  "

# Row 48 - FirstDigitalUSD
$ws.Range("E48").Value = "  +0.47%  "

# Row 49 - THORChain
Set-TextCell "D49" "8.51"
$ws.Range("E49").Value = "  -3.33%  "

# Row 50 - FLOKI
Set-TextCell "D50" "0.000247"
$ws.Range("E50").Value = "  +1.26%  "

# Row 51 - Mantle
Set-TextCell "D51" "1.30"
$ws.Range("E51").Value = "  +11.43%  "
